$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the ALU description strings (previously all used a[11..8]/b[15..12]
# which pointed at the wrong nibble of the operand; swap to a[15..12]/b[11..8]).
$ws.Range("D14").Value = "reg( a[15..12] ) <- a[11..0] / reg( a[11..8] ) + b[11..0] / reg( b[11..8] )"
$ws.Range("D15").Value = "reg( a[15..12] ) <- a[11..0] / reg( a[11..8] ) - b[11..0] / reg( b[11..8] )"
$ws.Range("D16").Value = "reg( a[15..12] ) <- a[11..0] / reg( a[11..8] ) * b[11..0] / reg( b[11..8] )"
$ws.Range("D17").Value = "reg( a[15..12] ) <- reg( a[11..8] ) >> b"

# D18 ("and" row) and D19 ("or" row) had their operator text swapped by mistake;
# fix by writing the correct operator text to each row.
$ws.Range("D18").Value = "reg( a[15..12] ) <- a[11..0] / reg( a[11..8] ) & b[11..0] / reg( b[11..8] )"
$ws.Range("D19").Value = "reg( a[15..12] ) <- a[11..0] / reg( a[11..8] ) | b[11..0] / reg( b[11..8] )"

$ws.Range("D20").Value = "reg( a[15..12] ) <- ! b[11..0] / !reg( b[15..12] )"

# Move the saved cursor/selection from D23 to D26.
$ws.Range("D26").Select()
